# Handles float input without breaking stuff
#
# This script fills in the previously-blank/placeholder marksheet with the
# student's actual results: summary counts (Right/Wrong/Not Attempt/Total),
# fixes the "-1" negative-marking value to be stored as a real number
# instead of text, records the student's per-question answers (with
# correct/incorrect/blank styling), and removes the unused 2nd/3rd
# "Student Ans / Correct Ans" blocks (columns D:E for rows 19+, and all of
# columns F:H) that were never filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# ---------------------------------------------------------------------
# 1) Summary block (rows 10-12): give the row labels in column A the
#    same bold "mtitleStyle" formatting already used by the header row
#    above them (row 9), and fill in the real tallies.
# ---------------------------------------------------------------------
$ws.Range("A9").Copy() | Out-Null
$excel.Union($ws.Range("A10"), $ws.Range("A11"), $ws.Range("A12")).PasteSpecial(-4122) | Out-Null

$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
# Negative marking was stored as text "-1"; make it a real number so
# downstream float/number handling doesn't break.
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "73/112"

# ---------------------------------------------------------------------
# 2) Per-question "Student Ans" column (A16:A40), compared against the
#    "Correct Ans" column B. Style mirrors the existing green/red/black
#    (correctStyle / incorrectStyle / normalStyle) cell styles already
#    used elsewhere on the sheet for this purpose.
# ---------------------------------------------------------------------
$excel.Union( `
    $ws.Range("A16"), $ws.Range("A18"), $ws.Range("A19"), $ws.Range("A21"), `
    $ws.Range("A23"), $ws.Range("A24"), $ws.Range("A25"), $ws.Range("A27"), `
    $ws.Range("A28"), $ws.Range("A29"), $ws.Range("A30"), $ws.Range("A32"), `
    $ws.Range("A33"), $ws.Range("A34"), $ws.Range("A35"), $ws.Range("A38"), `
    $ws.Range("A39") `
).Select() | Out-Null
$ws.Range("B10").Copy() | Out-Null
$excel.Selection.PasteSpecial(-4122) | Out-Null

$excel.Union($ws.Range("A22"), $ws.Range("A31")).Select() | Out-Null
$ws.Range("C10").Copy() | Out-Null
$excel.Selection.PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option A"
$ws.Range("A23").Value = "Option D"
$ws.Range("A24").Value = "Option A"
$ws.Range("A25").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option B"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A35").Value = "Option D"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
# Rows 17, 20, 26, 36, 37, 40 remain blank/not-attempted (already
# normalStyle with no value), so nothing to do for them.

# ---------------------------------------------------------------------
# 3) Second "Student Ans" column (D16:D18) - only the first three rows
#    of this block still have data; the rest (D19:E40) is removed below.
# ---------------------------------------------------------------------
$excel.Union($ws.Range("D16"), $ws.Range("D17")).Select() | Out-Null
$ws.Range("B10").Copy() | Out-Null
$excel.Selection.PasteSpecial(-4122) | Out-Null

$ws.Range("D18").Select() | Out-Null
$ws.Range("C10").Copy() | Out-Null
$excel.Selection.PasteSpecial(-4122) | Out-Null

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option B"

# ---------------------------------------------------------------------
# 4) Remove the now-unused parts of the 2nd block (D19:E40) and the
#    entire 3rd "Student Ans/Correct Ans" block (columns F:H), which
#    shrinks the sheet dimension down to A5:E40.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear() | Out-Null
$ws.Range("F1:H1048576").Clear() | Out-Null

$excel.CutCopyMode = 0
